$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update hours worked for week of row 10 (Thursday was 4.5, now 6.5; Friday added as 5)
$ws.Range("E10").Value = 6.5
$ws.Range("F10").Value = 5

# Recalculate so dependent formulas (I10 shared formula, I19 total) update
$excel.Calculate()

# Update the active selection on the sheet to O18
$ws.Range("O18").Select()
